$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column D entirely (this drops the "Ben" label in D1, shifts
# shared-string indices down, and shrinks the used range to A1:C26).
$ws.Range("D:D").Delete() | Out-Null

# Update the measured values in column C (rows 2-26) with the corrected,
# more precise readings.
$ws.Range("C2").Value = 10.5
$ws.Range("C3").Value = 10.5
$ws.Range("C4").Value = 10.5
$ws.Range("C5").Value = 10.5
$ws.Range("C6").Value = 10.5

$ws.Range("C7").Value = 51

$ws.Range("C8").Value = 50.9
$ws.Range("C9").Value = 50.9
$ws.Range("C10").Value = 50.9
$ws.Range("C11").Value = 50.9

$ws.Range("C12").Value = 100.2
$ws.Range("C13").Value = 100.2
$ws.Range("C14").Value = 100.2
$ws.Range("C15").Value = 100.2
$ws.Range("C16").Value = 100.3

$ws.Range("C17").Value = 201.2
$ws.Range("C18").Value = 201.2
$ws.Range("C19").Value = 201.2
$ws.Range("C20").Value = 201.2
$ws.Range("C21").Value = 201.2

$ws.Range("C22").Value = 298.4
$ws.Range("C23").Value = 298.4
$ws.Range("C24").Value = 298.3
$ws.Range("C25").Value = 298.3
$ws.Range("C26").Value = 298.3

# Update the current selection to mirror the saved state in the original file.
$ws.Range("N18").Select() | Out-Null
